# Rancangan Database Penjadwalan - "proccesing create data penjadwalan"
#
# Populates the new "Tb_Penjadwalan" table (columns AK:AO, rows 3-13) with
# scheduling data: No, Kode_Hari, Kode_Kelas, Kode_GMP, Kode_Sesi.
# Also re-orders the AL2:AO2 header row to match (Kode_Hari, Kode_Kelas,
# Kode_GMP, Kode_Sesi) and removes the old AN3:AN13 "jam mulai - jam selesai"
# helper formulas that are no longer needed on this range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): re-order Kode_Hari / Kode_Kelas / Kode_GMP / Kode_Sesi ---
$ws.Range("AL2").Value = "Kode_Hari"
$ws.Range("AM2").Value = "Kode_Kelas"
$ws.Range("AN2").Value = "Kode_GMP"
$ws.Range("AO2").Value = "Kode_Sesi"

# --- Data rows 3-13: No / Kode_Hari / Kode_Kelas / Kode_GMP / Kode_Sesi ---
# Clear the old time-range formulas first (AN3 unique, AN4:AN13 shared) -
# use Clear() (not ClearContents) so the leftover time-format style goes too.
$ws.Range("AN3:AN13").Clear()

$kodeSesi = @{
    3  = "sesi1"
    4  = "sesi2"
    5  = "sesi3"
    6  = "rest1"
    7  = "sesi4"
    8  = "sesi5"
    9  = "sesi6"
    10 = "rest2"
    11 = "sesi7"
    12 = "sesi8"
    13 = "sesi9"
}

for ($r = 3; $r -le 13; $r++) {
    $ws.Range("AK$r").Value = $r - 2
    $ws.Range("AL$r").Value = "h1"
    $ws.Range("AM$r").Value = "AP10"
    $ws.Range("AN$r").Value = "guru123"
    $ws.Range("AO$r").Value = $kodeSesi[$r]
}

# --- Sheet view: selection moved while entering the data ---
$ws.Range("AM4").Select()
